# Update "想去人数" (F column) values on the 展览 (sheet1) and 全部类型 (sheet4)
# worksheets to reflect newly generated output numbers.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 179
$ws1.Range("F5").Value = 3329
$ws1.Range("F6").Value = 341
$ws1.Range("F7").Value = 16
$ws1.Range("F8").Value = 421

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 179
$ws4.Range("F5").Value = 3329
$ws4.Range("F6").Value = 341
$ws4.Range("F9").Value = 16
$ws4.Range("F10").Value = 421
